{"js": "const replacements = [\n  [\"2025-09-20 Saturday\", \"2025-09-21 Sunday\"],\n  [\"111\u00d74=\", \"725\u00d78=\"],\n  [\"300\u00d74=\", \"357\u00d74=\"],\n  [\"633\u00d72=\", \"443\u00d74=\"],\n  [\"900\u00d78=\", \"899\u00d79=\"],\n  [\"441\u00d72=\", \"597\u00d73=\"],\n  [\"710\u00d78=\", \"569\u00d74=\"],\n  [\"635\u00d76=\", \"239\u00d75=\"],\n  [\"407\u00d78=\", \"274\u00d76=\"],\n  [\"807\u00d74=\", \"498\u00d79=\"],\n  [\"734\u00d74=\", \"960\u00d77=\"],\n  [\"443\u00d76=\", \"669\u00d72=\"],\n  [\"687\u00d78=\", \"205\u00d74=\"],\n  [\"977\u00d76=\", \"372\u00d78=\"],\n  [\"817\u00d78=\", \"443\u00d77=\"],\n  [\"397\u00d75=\", \"163\u00d79=\"],\n  [\"484\u00d78=\", \"329\u00d78=\"],\n  [\"333\u00d79=\", \"580\u00d73=\"],\n  [\"739\u00d74=\", \"406\u00d72=\"],\n  [\"340\u00d75=\", \"535\u00d76=\"],\n  [\"471\u00d77=\", \"666\u00d75=\"],\n  [\"455\u00d75=\", \"289\u00d73=\"],\n  [\"176\u00d76=\", \"989\u00d76=\"],\n  [\"362\u00d73=\", \"868\u00d75=\"],\n  [\"416\u00d77=\", \"490\u00d75=\"],\n  [\"627\u00d76=\", \"495\u00d74=\"],\n];\n\nfor (const [oldText, newText] of replacements) {\n  const results = context.document.body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n  for (const item of results.items) {\n    item.insertText(newText, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n", "ps1": "$d = $word.ActiveDocument\n\nfunction Replace-Text($old, $new) {\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Text = $old\n    $find.Replacement.Text = $new\n    $find.Execute($find.Text, $false, $false, $false, $false, $false, $true, 1, $false, $find.Replacement.Text, 2)\n}\n\nReplace-Text '2025-09-20 Saturday' '2025-09-21 Sunday'\nReplace-Text '111\u00d74=' '725\u00d78='\nReplace-Text '300\u00d74=' '357\u00d74='\nReplace-Text '633\u00d72=' '443\u00d74='\nReplace-Text '900\u00d78=' '899\u00d79='\nReplace-Text '441\u00d72=' '597\u00d73='\nReplace-Text '710\u00d78=' '569\u00d74='\nReplace-Text '635\u00d76=' '239\u00d75='\nReplace-Text '407\u00d78=' '274\u00d76='\nReplace-Text '807\u00d74=' '498\u00d79='\nReplace-Text '734\u00d74=' '960\u00d77='\nReplace-Text '443\u00d76=' '669\u00d72='\nReplace-Text '687\u00d78=' '205\u00d74='\nReplace-Text '977\u00d76=' '372\u00d78='\nReplace-Text '817\u00d78=' '443\u00d77='\nReplace-Text '397\u00d75=' '163\u00d79='\nReplace-Text '484\u00d78=' '329\u00d78='\nReplace-Text '333\u00d79=' '580\u00d73='\nReplace-Text '739\u00d74=' '406\u00d72='\nReplace-Text '340\u00d75=' '535\u00d76='\nReplace-Text '471\u00d77=' '666\u00d75='\nReplace-Text '455\u00d75=' '289\u00d73='\nReplace-Text '176\u00d76=' '989\u00d76='\nReplace-Text '362\u00d73=' '868\u00d75='\nReplace-Text '416\u00d77=' '490\u00d75='\nReplace-Text '627\u00d76=' '495\u00d74='\n"}
